$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking strings
# (e.g. "0.9998") are stored as text, matching the workbook's inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.484.48"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.913.66"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "244.63"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4831"
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("D8").Value = "0.2890"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "0.06720"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").Value = "109.77"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").Value = "19.06"
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("D12").Value = "1.912.86"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "0.07550"
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "5.267"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "0.6723"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").Value = "283.30"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "30.489.37"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "0.000007574"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "2.168.19"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "5.476"
$ws.Range("E22").Value = "  +4.67%  "
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "6.432"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "9.447"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "164.45"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "20.21"
$ws.Range("E27").Value = "  -7.08%  "
$ws.Range("D28").Value = "2.119"
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "1.405"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").Value = "4.166"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "4.037"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("D33").Value = "0.04987"
$ws.Range("E33").Value = "  -1.41%  "
$ws.Range("D34").Value = "0.7314"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "0.02029"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "110.62"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "2.018"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "0.4456"
$ws.Range("E42").Value = "  +5.37%  "
$ws.Range("D43").Value = "0.8656"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "5.792"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "68.03"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "7.329"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "49.07"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").Value = "9.288"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "0.1241"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "34.81"
$ws.Range("E51").Value = "  -0.06%  "

# Restore the original (default) style on column D so no stray number-format
# style is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
